$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.420.63'
$ws.Range('E2').Value = '  +0.74%  '
$ws.Range('D3').Value = '1.920.77'
$ws.Range('E3').Value = '  +3.44%  '
$ws.Range('E4').Value = '  -0.46%  '
$ws.Range('D5').Value = "'244.55"
$ws.Range('E5').Value = '  +2.33%  '
$ws.Range('D6').Value = "'0.659"
$ws.Range('E6').Value = '  +5.84%  '
$ws.Range('E7').Value = '  -0.42%  '
$ws.Range('D8').Value = "'41.62"
$ws.Range('E8').Value = '  -0.48%  '
$ws.Range('E9').Value = '  +7.13%  '
$ws.Range('D10').Value = "'52.82"
$ws.Range('E10').Value = '  +12.58%  '
$ws.Range('D11').Value = "'0.0717"
$ws.Range('E11').Value = '  +3.58%  '
$ws.Range('D12').Value = "'0.0997"
$ws.Range('E12').Value = '  +0.85%  '
$ws.Range('D13').Value = '2.195.29'
$ws.Range('E13').Value = '  +3.26%  '
$ws.Range('D14').Value = "'12.10"
$ws.Range('E14').Value = '  +5.62%  '
$ws.Range('D15').Value = "'0.701"
$ws.Range('E15').Value = '  +3.81%  '
$ws.Range('D16').Value = '1.926.16'
$ws.Range('E16').Value = '  +3.77%  '
$ws.Range('D17').Value = "'4.88"
$ws.Range('E17').Value = '  +3.65%  '
$ws.Range('D18').Value = '35.425.05'
$ws.Range('E18').Value = '  +0.82%  '
$ws.Range('D19').Value = "'72.22"
$ws.Range('E19').Value = '  +3.67%  '
$ws.Range('D20').Value = '0.0₃0825'
$ws.Range('E20').Value = '  +3.77%  '
$ws.Range('D21').Value = "'240.37"
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('E22').Value = '  +2.51%  '
$ws.Range('D23').Value = "'4.89"
$ws.Range('E23').Value = '  +3.24%  '
$ws.Range('D25').Value = "'2.29"
$ws.Range('E25').Value = '  +1.29%  '
$ws.Range('D26').Value = "'2.31"
$ws.Range('E26').Value = '  +20.85%  '
$ws.Range('D27').Value = "'170.64"
$ws.Range('E27').Value = '  +1.08%  '
$ws.Range('D28').Value = "'8.46"
$ws.Range('E28').Value = '  +6.17%  '
$ws.Range('D29').Value = "'18.54"
$ws.Range('E29').Value = '  +5.05%  '
$ws.Range('E30').Value = '  +2.37%  '
$ws.Range('D31').Value = "'4.14"
$ws.Range('E31').Value = '  +3.80%  '
$ws.Range('E32').Value = '  +1.78%  '
$ws.Range('D33').Value = "'0.941"
$ws.Range('E33').Value = '  +13.41%  '
$ws.Range('E34').Value = '  -0.43%  '
$ws.Range('E35').Value = '  +3.18%  '
$ws.Range('D36').Value = "'1.75"
$ws.Range('E36').Value = '  -5.03%  '
$ws.Range('D37').Value = "'2.06"
$ws.Range('E37').Value = '  +2.86%  '
$ws.Range('E38').Value = '  +1.66%  '
$ws.Range('D39').Value = "'1.12"
$ws.Range('E39').Value = '  +2.94%  '
$ws.Range('E41').Value = '  +17.14%  '
$ws.Range('D42').Value = "'16.27"
$ws.Range('E42').Value = '  +8.97%  '
$ws.Range('D43').Value = "'90.81"
$ws.Range('E43').Value = '  +1.34%  '
$ws.Range('D44').Value = '1.345.91'
$ws.Range('E44').Value = '  +0.33%  '
$ws.Range('E45').Value = '  +2.98%  '
$ws.Range('D46').Value = "'48.73"
$ws.Range('E46').Value = '  +40.09%  '
$ws.Range('E47').Value = '  +2.62%  '
$ws.Range('E48').Value = '  -0.07%  '
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('D50').Value = '2.106.24'
$ws.Range('E50').Value = '  +3.31%  '
$ws.Range('D51').Value = "'0.0700"
$ws.Range('E51').Value = '  +3.04%  '
